$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 41 (pushes the existing rows 41..136 down to 42..137).
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new weekly record.
$ws.Range("A41").Value2 = 7
$ws.Range("B41").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C41").Value2 = "Ñuble"
$ws.Range("D41").Value2 = 45070
$ws.Range("E41").Value2 = 16
$ws.Range("F41").Value2 = 100112031
$ws.Range("G41").Value2 = "Poroto verde"
$ws.Range("H41").Value2 = "Magnum"
$ws.Range("I41").Value2 = "Primera"
$ws.Range("J41").Value2 = 20
$ws.Range("K41").Value2 = 35000
$ws.Range("L41").Value2 = 35000
$ws.Range("M41").Value2 = 35000
$ws.Range("N41").Value2 = "`$/malla 25 kilos"
$ws.Range("O41").Value2 = "Provincia de Limarí"
$ws.Range("P41").Value2 = 1400
$ws.Range("Q41").Value2 = 25
$ws.Range("R41").Value2 = "Hortaliza"
